$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: rename interview, swap Pass/Fail counts
$ws.Range("B3").Value = "1725397919-RAG-Average"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 10

# Row 6: rename interview, update Pass/Fail counts
$ws.Range("B6").Value = "1725380262-Receptionist-Average"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 8

# Update the selected cell to match the saved view state
$ws.Range("C3").Select()
